$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Row 7 / column B ("Experimental" value) was empty; it should now read
# the literal text "false". Assigning the raw string directly would be
# auto-typed to a Boolean by the engine (mirrors real Excel "TRUE"/"FALSE"
# literal coercion), so route it through a text-producing formula and
# collapse it back down to a static value via Copy/PasteSpecial (values
# only) - this keeps the cell's string type and preserves its existing
# style (s="2") instead of minting a new quote-prefixed style.
$b7 = $ws1.Range("B7")
$b7.Formula = '=""&"false"'
$b7.Copy() | Out-Null
$b7.PasteSpecial(-4163) | Out-Null   # xlPasteValues

# Row 8 / column B: bump the recorded Date property to the new timestamp.
$ws1.Range("B8").Value = "2025-11-30T13:08:37+00:00"
